$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44334
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12500
$ws.Range("Q2").Value = "$/caja 12 kilos empedrada"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1042
$ws.Range("T2").Value = 12

# Row 4
$ws.Range("D4").Value = 44330
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 15500
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 861
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("S5").Value = 1109

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = "$/caja 16 kilos granel"
$ws.Range("S6").Value = 1000
$ws.Range("T6").Value = 16
